# HIVE TEAMS.docx -> Dutch translation pass
# Applies the set of English -> Dutch text replacements described by the
# commit "New translations HIVE TEAMS.docx (Dutch)".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Developer" -> "Ontwikkelaar": this exact paragraph text occurs 8
#    times in the document, but only the 4 instances inside the
#    "HIVE TEAM: WEB" section (before "HIVE TEAM: QUALITY ASSURANCE")
#    are translated; the 4 instances later in "HIVE TEAM: DEVELOPMENT"
#    are left in English. Walk paragraphs and only touch the ones in
#    that scoped section. This runs first, while the section heading
#    text is still in English, so the section boundary check is simple
#    and unambiguous.
# ---------------------------------------------------------------------

$inWebTeam = $false
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "HIVE TEAM: WEB") {
        $inWebTeam = $true
    } elseif ($ptext -like "*HIVE TEAM: QUALITY ASSURANCE*") {
        $inWebTeam = $false
    } elseif ($inWebTeam -and $ptext -eq "Developer") {
        $p.Range.Text = "Ontwikkelaar"
    }
}

# ---------------------------------------------------------------------
# 2) Simple, uniquely-occurring paragraph/sentence replacements.
#    Each old string appears exactly once in the document, so a direct
#    Find/Replace on the whole story is safe. We set Range.Text (rather
#    than using Find's built-in Replace) so Word's "smart quotes"
#    autocorrect doesn't mangle the straight apostrophes/quotes that the
#    Dutch copy intentionally uses.
# ---------------------------------------------------------------------

function Replace-UniqueText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $newText
    } else {
        Write-Output "NOT FOUND: $oldText"
    }
    return $found
}

Replace-UniqueText `
    "Van de Getto's van Brazilië tot de cryptoverse. Father, husband and technolover, born in Brazil and living in Australia. Mobile and Web Development." `
    "Van de Getto's van Brazilië tot de cryptoverse. Vader, echtgenoot en techno-liefhebber, geboren in Brazilië en woonachtig in Australië. Mobiele en webontwikkeling."

Replace-UniqueText "Front-End Developer" "Front-end ontwikkelaar"

Replace-UniqueText "UX/Product Design" "UX/productontwerp"

Replace-UniqueText "HIVE TEAM: QUALITY ASSURANCE" "HIVE TEAM: KWALITEITSBORGING"

Replace-UniqueText "Ensuring all development tasks meet quality criteria." "Zorgen dat alle ontwikkeltaken voldoen aan kwaliteitscriteria."

Replace-UniqueText "WANNA GET INVOLVED?" "WIL JE MEEDOEN?"

Replace-UniqueText "The SmartHive has a place for folks of all backgrounds. Come hungry!" "De SmartHive heeft een plek voor mensen van alle achtergronden. Neem je passie mee!"

Replace-UniqueText `
    ("We believe " + [char]0x2018 + "Core" + [char]0x2019 + " teams are a bad idea and something that ultimately leads to inefficiency and corruption. We want to move past it and create a decentralized organizational model inspired by ant and bee colonies.") `
    "Wij geloven dat 'Core'-teams een slecht idee zijn en uiteindelijk leiden tot inefficiëntie en corruptie. We willen dat achter ons laten en een gedecentraliseerd organisatiemodel creëren dat geïnspireerd is op mieren- en bijenkolonies."

Replace-UniqueText `
    "In order to create and maintain a decentralized governance structure, we are introducing two concepts SmartHive and Hive Structuring Teams (HST). SmartHive enables anyone that holds coins the opportunity to vote on proposals submitted by the community. SmartHive will be the lifeblood of the project, which will allow anyone to get involved and submit proposals, helping to generate organic growth at a grassroots level, creating a bottom-up management structure." `
    ("Om een " + [char]0x200B + [char]0x200B + "gedecentraliseerde bestuursstructuur te creëren en te behouden, introduceren we twee concepten, SmartHive en Hive Structuring Teams (HST). SmartHive stelt iedereen die munten bezit in staat om te stemmen over voorstellen die door de gemeenschap zijn ingediend. SmartHive will be the lifeblood of the project, which will allow anyone to get involved and submit proposals, helping to generate organic growth at a grassroots level, creating a bottom-up management structure.")

# ---------------------------------------------------------------------
# 3) Drop the stray leading non-breaking-space run that used to sit in
#    front of the "HIVE TEAM: QUALITY ASSURANCE" heading run (that
#    paragraph is the only heading with the nbsp split into its own
#    run instead of being merged into the heading text). Match on the
#    already-translated Dutch heading since this step runs after (2).
# ---------------------------------------------------------------------

foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext -like "*HIVE TEAM: KWALITEITSBORGING*") {
        $pr = $p.Range
        $found = $pr.Find.Execute([char]0x00A0, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $pr.Delete()
        }
    }
}
